$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has one data table, rows 2..263 (row 1 = header). Two brand-new
# data rows are being inserted right after the current row 163, pushing the
# old rows 164..263 down to 166..265. Two more brand-new rows are appended
# at the very end (new rows 264 and 265).

# --- 1) Insert two blank rows before row 164 (shifts old 164..263 -> 166..265)
$ws.Rows.Item(164).Resize(2).Insert()

# --- 2) Fill the first new row (164)
$ws.Cells.Item(164,1).Value = 5
$ws.Cells.Item(164,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(164,3).Value = "Maule"
$ws.Cells.Item(164,4).Value = 44452
$ws.Cells.Item(164,5).Value = 7
$ws.Cells.Item(164,6).Value = 100114001
$ws.Cells.Item(164,7).Value = "Papa"
$ws.Cells.Item(164,8).Value = "Asterix"
$ws.Cells.Item(164,9).Value = "1a (guarda)"
$ws.Cells.Item(164,10).Value = 1200
$ws.Cells.Item(164,11).Value = 8000
$ws.Cells.Item(164,12).Value = 8000
$ws.Cells.Item(164,13).Value = 8000
$ws.Cells.Item(164,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(164,15).Value = "Región del Maule"
$ws.Cells.Item(164,16).Value = 320
$ws.Cells.Item(164,17).Value = 25
$ws.Cells.Item(164,18).Value = "Hortaliza"

# --- 3) Fill the second new row (165)
$ws.Cells.Item(165,1).Value = 5
$ws.Cells.Item(165,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(165,3).Value = "Maule"
$ws.Cells.Item(165,4).Value = 44452
$ws.Cells.Item(165,5).Value = 7
$ws.Cells.Item(165,6).Value = 100114001
$ws.Cells.Item(165,7).Value = "Papa"
$ws.Cells.Item(165,8).Value = "Rodeo"
$ws.Cells.Item(165,9).Value = "1a (guarda lavada)"
$ws.Cells.Item(165,10).Value = 1200
$ws.Cells.Item(165,11).Value = 10000
$ws.Cells.Item(165,12).Value = 10000
$ws.Cells.Item(165,13).Value = 10000
$ws.Cells.Item(165,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(165,15).Value = "Región de Los Lagos"
$ws.Cells.Item(165,16).Value = 400
$ws.Cells.Item(165,17).Value = 25
$ws.Cells.Item(165,18).Value = "Hortaliza"

# --- 4) Append two brand-new rows at the end (264 and 265)
$ws.Cells.Item(264,1).Value = 5
$ws.Cells.Item(264,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(264,3).Value = "Maule"
$ws.Cells.Item(264,4).Value = 44399
$ws.Cells.Item(264,5).Value = 7
$ws.Cells.Item(264,6).Value = 100114001
$ws.Cells.Item(264,7).Value = "Papa"
$ws.Cells.Item(264,8).Value = "Rodeo"
$ws.Cells.Item(264,9).Value = "1a (guarda lavada)"
$ws.Cells.Item(264,10).Value = 1200
$ws.Cells.Item(264,11).Value = 8000
$ws.Cells.Item(264,12).Value = 8000
$ws.Cells.Item(264,13).Value = 8000
$ws.Cells.Item(264,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(264,15).Value = "Región de Los Lagos"
$ws.Cells.Item(264,16).Value = 320
$ws.Cells.Item(264,17).Value = 25
$ws.Cells.Item(264,18).Value = "Hortaliza"
$ws.Cells.Item(264,4).NumberFormat = $ws.Cells.Item(263,4).NumberFormat

$ws.Cells.Item(265,1).Value = 5
$ws.Cells.Item(265,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(265,3).Value = "Maule"
$ws.Cells.Item(265,4).Value = 44400
$ws.Cells.Item(265,5).Value = 7
$ws.Cells.Item(265,6).Value = 100114001
$ws.Cells.Item(265,7).Value = "Papa"
$ws.Cells.Item(265,8).Value = "Rodeo"
$ws.Cells.Item(265,9).Value = "1a (guarda lavada)"
$ws.Cells.Item(265,10).Value = 1200
$ws.Cells.Item(265,11).Value = 8000
$ws.Cells.Item(265,12).Value = 8000
$ws.Cells.Item(265,13).Value = 8000
$ws.Cells.Item(265,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(265,15).Value = "Región de La Araucanía"
$ws.Cells.Item(265,16).Value = 320
$ws.Cells.Item(265,17).Value = 25
$ws.Cells.Item(265,18).Value = "Hortaliza"
$ws.Cells.Item(265,4).NumberFormat = $ws.Cells.Item(263,4).NumberFormat
